$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.128.38"

$ws.Range("E2").Value = "  -1.85%  "

$ws.Range("D3").Value = "2.338.22"

$ws.Range("E3").Value = "  -4.53%  "

$ws.Range("E4").Value = "  +0.23%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "540.32"
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").Value = "  -0.89%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.01"
$ws.Range("D6").Style = "Normal"

$ws.Range("E6").Value = "  -6.64%  "

$ws.Range("E7").Value = "  +0.17%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.519"
$ws.Range("D8").Style = "Normal"

$ws.Range("E8").Value = "  -11.21%  "

$ws.Range("D9").Value = "2.335.84"

$ws.Range("E9").Value = "  -4.59%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.103"
$ws.Range("D10").Style = "Normal"

$ws.Range("E10").Value = "  -2.56%  "

$ws.Range("E11").Value = "  -0.23%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.24"
$ws.Range("D12").Style = "Normal"

$ws.Range("E12").Value = "  -2.89%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.338"
$ws.Range("D13").Style = "Normal"

$ws.Range("E13").Value = "  -3.57%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "24.27"
$ws.Range("D14").Style = "Normal"

$ws.Range("E14").Value = "  -6.36%  "

$ws.Range("D15").Value = "2.768.48"

$ws.Range("E15").Value = "  -3.86%  "

$ws.Range("D16").Value = "60.113.12"

$ws.Range("E16").Value = "  -1.64%  "

$ws.Range("E17").Value = "  -3.90%  "

$ws.Range("D18").Value = "2.345.64"

$ws.Range("E18").Value = "  -5.32%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.46"
$ws.Range("D19").Style = "Normal"

$ws.Range("E19").Value = "  -4.63%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.03"
$ws.Range("D20").Style = "Normal"

$ws.Range("E20").Value = "  -2.97%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "309.40"
$ws.Range("D21").Style = "Normal"

$ws.Range("E21").Value = "  -2.68%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.52"
$ws.Range("D22").Style = "Normal"

$ws.Range("E22").Value = "  -6.35%  "

$ws.Range("E23").Value = "  +0.04%  "

$ws.Range("E24").Value = "  -1.37%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "62.73"
$ws.Range("D25").Style = "Normal"

$ws.Range("E25").Value = "  -1.16%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.47"
$ws.Range("D26").Style = "Normal"

$ws.Range("E26").Value = "  +8.63%  "

$ws.Range("E27").Value = "  +0.14%  "

$ws.Range("D28").Value = "2.462.00"

$ws.Range("E28").Value = "  -3.51%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.87"
$ws.Range("D29").Style = "Normal"

$ws.Range("E29").Value = "  -4.22%  "

$ws.Range("D30").Value = "0.0₃0872"

$ws.Range("E30").Value = "  -10.03%  "

$ws.Range("E31").Value = "  -6.49%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "490.89"
$ws.Range("D32").Style = "Normal"

$ws.Range("E32").Value = "  -6.86%  "

$ws.Range("E33").Value = "  -2.95%  "

$ws.Range("E34").Value = "  -5.23%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.50"
$ws.Range("D35").Style = "Normal"

$ws.Range("E35").Value = "  -4.39%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.998"
$ws.Range("D36").Style = "Normal"

$ws.Range("E36").Value = "  -0.05%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.52"
$ws.Range("D37").Style = "Normal"

$ws.Range("E37").Value = "  -5.31%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.370"
$ws.Range("D38").Style = "Normal"

$ws.Range("E38").Value = "  -1.99%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.15"
$ws.Range("D39").Style = "Normal"

$ws.Range("E39").Value = "  -0.26%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.18"
$ws.Range("D40").Style = "Normal"

$ws.Range("E40").Value = "  -9.25%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.76"
$ws.Range("D41").Style = "Normal"

$ws.Range("E41").Value = "  +0.47%  "

$ws.Range("E42").Value = "  -0.03%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "136.27"
$ws.Range("D43").Style = "Normal"

$ws.Range("E43").Value = "  -4.26%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "39.90"
$ws.Range("D44").Style = "Normal"

$ws.Range("E44").Value = "  -0.51%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "140.69"
$ws.Range("D45").Style = "Normal"

$ws.Range("E45").Value = "  -0.33%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.08"
$ws.Range("D46").Style = "Normal"

$ws.Range("E46").Value = "  -8.67%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.50"
$ws.Range("D47").Style = "Normal"

$ws.Range("E47").Value = "  -2.86%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0505"
$ws.Range("D48").Style = "Normal"

$ws.Range("E48").Value = "  -5.40%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "19.23"
$ws.Range("D49").Style = "Normal"

$ws.Range("E49").Value = "  -9.24%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.565"
$ws.Range("D50").Style = "Normal"

$ws.Range("E50").Value = "  -3.36%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0890"
$ws.Range("D51").Style = "Normal"

$ws.Range("E51").Value = "  -4.37%  "
